$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @("Configuration setup", "Read all the configurations", "Pass", "22-01-2021 10:31:03 PM"),
    @("Configuration setup", "Read all the configurations", "Pass", "25-01-2021 11:45:30 AM"),
    @("Login to Portal ", "Logged in sucessfully ", "Pass", "25-01-2021 11:46:04 AM"),
    @("Navigation to LATAM", "Successfully Navigated to LATAM ", "Pass", "25-01-2021 11:46:19 AM"),
    @("Search of created WO ", "Search is Working", "Pass", "25-01-2021 11:46:31 AM"),
    @("Segment level Ingest", "Segment level Ingest is completed", "Pass", "25-01-2021 11:46:59 AM"),
    @("AQC overRide is clicked ", "AQC over Ride performed ", "Pass", "25-01-2021 11:47:14 AM"),
    @("AQC overRide pass is enabled ", "AQC over ride Pass", "Fail", "25-01-2021 11:47:20 AM"),
    @("Configuration setup", "Read all the configurations", "Pass", "25-01-2021 12:42:13 PM"),
    @("Login to Portal ", "Logged in sucessfully ", "Pass", "25-01-2021 12:42:43 PM"),
    @("Configuration setup", "Read all the configurations", "Pass", "25-01-2021 12:44:28 PM"),
    @("Login to Portal ", "Logged in sucessfully ", "Pass", "25-01-2021 12:44:59 PM"),
    @("Navigation to LATAM", "Successfully Navigated to LATAM ", "Pass", "25-01-2021 12:45:14 PM"),
    @("Search of created WO ", "Search is Working", "Pass", "25-01-2021 12:45:27 PM"),
    @("Segment level Ingest", "Segment level Ingest is completed", "Pass", "25-01-2021 12:45:54 PM"),
    @("AQC overRide is clicked ", "AQC over Ride performed ", "Pass", "25-01-2021 12:46:13 PM"),
    @("AQC overRide Report Download ", "AQC overRide Report Download ", "Pass", "25-01-2021 12:46:20 PM"),
    @("AQC overRide pass is enabled ", "AQC over ride Pass", "Fail", "25-01-2021 12:46:20 PM")
)

$startRow = 141
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $row = $rows[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
}
